# Update the "Confirmation Number" (column F) values in the
# "Form responses 1" sheet with the new, de-duplicated ID numbers
# produced after re-running the RPA workflow / integrating with the
# team's datatables.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form responses 1")

$newValues = @(
    14583376,
    15480124,
    19760685,
    22789443,
    24219363,
    24770454,
    26663693,
    27498330,
    28557906,
    28688027,
    30279019,
    30338111,
    30990800,
    31310453,
    34344874,
    34650095,
    36373341,
    37105247,
    38559467,
    39467735,
    41370596,
    42487189,
    42910376,
    45345003,
    46587488,
    46935167,
    53195057,
    53416701,
    53596124,
    54076152,
    58574868,
    62193437,
    62283126,
    62798342,
    67424150,
    67431110,
    69526335,
    71034050,
    72696126,
    74719844,
    78475997,
    78668846,
    80699407,
    81218466,
    86347520,
    87857672,
    88062606,
    89431443,
    91762639,
    92312909,
    93038340,
    93268593,
    94192223,
    95425054,
    98427834,
    99338799
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 6).Value = $newValues[$i]
}
